$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 00:27"

# Row 4
$ws.Range("B4").Value = 2497544
$ws.Range("C4").Value = 34990
$ws.Range("D4").Value = 1047041
$ws.Range("E4").Value = 1325654
$ws.Range("G4").Value = 568
$ws.Range("H4").Value = 124849

# Row 5
$ws.Range("B5").Value = 1228114
$ws.Range("C5").Value = 35640
$ws.Range("E5").Value = 523235
$ws.Range("G5").Value = 1097
$ws.Range("H5").Value = 54971

# Row 15
$ws.Range("B15").Value = 193715
$ws.Range("C15").Value = 461
$ws.Range("E15").Value = 7903

# Row 22
$ws.Range("B22").Value = 102576
$ws.Range("C22").Value = 334
$ws.Range("D22").Value = 65367
$ws.Range("E22").Value = 28708

# Row 25
$ws.Range("B25").Value = 80599
$ws.Range("C25").Value = 3486
$ws.Range("D25").Value = 33349
$ws.Range("E25").Value = 44596
$ws.Range("G25").Value = 163
$ws.Range("H25").Value = 2654

# Row 55
$ws.Range("B55").Value = 18110
$ws.Range("C55").Value = 86
$ws.Range("D55").Value = 16320
$ws.Range("E55").Value = 822
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 968

# Row 88
$ws.Range("A88").Value = "Bulgaria"
$ws.Range("B88").Value = 4408
$ws.Range("C88").Value = 166
$ws.Range("D88").Value = 2370
$ws.Range("E88").Value = 1827
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 211

# Row 89
$ws.Range("A89").Value = "Venezuela"
$ws.Range("B89").Value = 4366
$ws.Range("D89").Value = 1327
$ws.Range("E89").Value = 3001
$ws.Range("H89").Value = 38

# Row 107
$ws.Range("A107").Value = "Mali"
$ws.Range("B107").Value = 2039
$ws.Range("C107").Value = 34
$ws.Range("D107").Value = 1383
$ws.Range("E107").Value = 543
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 113

# Row 108
$ws.Range("A108").Value = "Sri Lanka"
$ws.Range("B108").Value = 2010
$ws.Range("C108").Value = 9
$ws.Range("D108").Value = 1602
$ws.Range("E108").Value = 397
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 11

# Row 168
$ws.Range("D168").Value = 184
$ws.Range("E168").Value = 11

# Row 202
$ws.Range("A202").Value = "Dominica"

# Row 203
$ws.Range("A203").Value = "Fiyi"

# Row 208
$ws.Range("A208").Value = "Islas Malvinas"

# Row 209
$ws.Range("A209").Value = "Groenlandia"

# Row 211
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 212
$ws.Range("A212").Value = "Seychelles"
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

